# Insert a new price-report row for "Ajo" (Chino / Primera) at row 456.
# Excel shifts every existing row at/after 456 down by one (456->457, ...,
# 566->567), which is exactly what the target diff shows: the sheet's used
# range grows from A1:R566 to A1:R567 and every row from 457 on carries the
# values that used to sit one row above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(456).Insert()

$ws.Range("A456").Value = 5
$ws.Range("B456").Value = 'Macroferia Regional de Talca'
$ws.Range("C456").Value = 'Maule'
$ws.Range("D456").Value = 45244
$ws.Range("E456").Value = 7
$ws.Range("F456").Value = 100112003
$ws.Range("G456").Value = 'Ajo'
$ws.Range("H456").Value = 'Chino'
$ws.Range("I456").Value = 'Primera'
$ws.Range("J456").Value = 200
$ws.Range("K456").Value = 23000
$ws.Range("L456").Value = 23000
$ws.Range("M456").Value = 23000
$ws.Range("N456").Value = '$/malla 10 kilos'
$ws.Range("O456").Value = 'China'
$ws.Range("P456").Value = 2300
$ws.Range("Q456").Value = 10
$ws.Range("R456").Value = 'Hortaliza'
